# Generate Report for Archive
#
# The handoff status text changes from "Ready for handoff" to
# "In Translation" everywhere it appears (Overview!E2:F4, zh-cn!C2:C4,
# de-de!C2:C4 all share the same string). Replacing the text shrinks the
# column(s) that hold it, so the affected "Status" columns are re-sized
# to match the new, shorter content.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# New column width (Excel "characters" units) that the shortened text
# autofits down to: 13.4101845877511 stored width - the 5/6 character
# fixed padding Excel always adds on top of ColumnWidth.
$newColumnWidth = 13.4101845877511 - 0.8333333333333334

foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace($oldStatus, $newStatus)
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1:F1").EntireColumn.ColumnWidth = $newColumnWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").EntireColumn.ColumnWidth = $newColumnWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").EntireColumn.ColumnWidth = $newColumnWidth
